$d = $word.ActiveDocument

# Old/new text pairs, in the order the cells appear in the document.
$replacements = @(
    "96÷5=19, 1", "39÷8=4, 7",
    "44÷7=6, 2", "75÷7=10, 5",
    "62÷2=31, 0", "96÷8=12, 0",
    "49÷3=16, 1", "12÷8=1, 4",
    "79÷4=19, 3", "17÷3=5, 2",
    "49÷4=12, 1", "58÷4=14, 2",
    "57÷3=19, 0", "29÷4=7, 1",
    "10÷8=1, 2", "92÷8=11, 4",
    "83÷3=27, 2", "22÷9=2, 4",
    "75÷8=9, 3", "72÷9=8, 0",
    "30÷6=5, 0", "38÷5=7, 3",
    "94÷4=23, 2", "23÷8=2, 7",
    "91÷4=22, 3", "48÷5=9, 3",
    "64÷2=32, 0", "67÷5=13, 2",
    "70÷5=14, 0", "19÷3=6, 1",
    "53÷7=7, 4", "52÷8=6, 4",
    "22÷3=7, 1", "55÷6=9, 1",
    "23÷8=2, 7", "65÷4=16, 1",
    "45÷5=9, 0", "91÷6=15, 1",
    "22÷3=7, 1", "19÷5=3, 4",
    "35÷3=11, 2", "10÷2=5, 0",
    "44÷9=4, 8", "81÷7=11, 4",
    "82÷9=9, 1", "60÷7=8, 4",
    "47÷6=7, 5", "91÷8=11, 3",
    "55÷5=11, 0", "40÷4=10, 0"
)

$table = $d.Tables.Item(1)
$idx = 0
foreach ($row in $table.Rows) {
    foreach ($cell in $row.Cells) {
        $cellRange = $cell.Range
        $cellText = $cellRange.Text
        # cell text includes trailing cell-mark characters; trim them
        $trimmed = $cellText.TrimEnd([char]0x0007, [char]0x000D)
        if ($trimmed.Length -gt 0) {
            $old = $replacements[$idx]
            $new = $replacements[$idx + 1]
            if ($trimmed -ne $old) {
                throw "Mismatch at cell index $idx : expected '$old' but found '$trimmed'"
            }
            # Scope the Find to exactly this cell's character range, and use
            # wdReplaceOne (1) rather than wdReplaceAll (2) so the
            # replacement touches only the single match found inside this
            # cell -- never other cells that may (at this moment, or after
            # this edit) contain identical text elsewhere in the document.
            $scoped = $d.Range($cellRange.Start, $cellRange.End)
            $scoped.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1)
            $idx += 2
        }
    }
}
